$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Article 96 (Chapter 4 / "5 Apr 2020") now has a finished/live session log:
#  - the existing Kick Start speaker line gets a trailing "Download" action label
#  - a new speaker line is added for the second talk of the day (Mona Majeed / Prezi link)
$text = "Date; 5 Apr 2020`r`n" +
    "Title; Chapter 4: Open Mic Morning in Pakistan`r`n" +
    "Video; Not Available yet`r`n" +
    "Duration; 0314 hrs`r`n" +
    "Speaker; 1030 hrs, Kick Start by Qasim Ali, techshek4.pptx, Download`r`n" +
    "Speaker; 1400 hrs, Dedicated Parents by Mona Majeed, https://prezi.com/view/xtik0a6jQUqf2BIcVjio/, Open Link`r`n" +
    "Pictures; Not Found"

$cell = $ws.Range("A4")
$cell.Value = $text

# The cell grew from six to seven wrapped lines, so the row needs to grow
# from its previous 136pt to fit the extra "Speaker" line (204pt).
$ws.Rows.Item(4).RowHeight = 204
